# Update the 2022-07-13 (row 15) daily figures for Lanzhou covid-19 data.
# Source data for that day had just come in: 22 new confirmed, 130 new
# asymptomatic, 20 Chengguan confirmed, 121 Chengguan asymptomatic.
# The C31:F31 totals are SUM() formulas over C5:C29/D5:D29/E5:E29/F5:F29,
# so they recalculate automatically once the inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C15").Value = 22
$ws.Range("D15").Value = 130
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 121

# Leave the cursor where the editor last left it.
$ws.Range("J16").Select()
